# Automatische test-sync: 2025-08-28 21:24:50
# Appends a new log row to "Logs" and updates the "Dashboard" summary table.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 30 to the Logs sheet ---
$newRow = 30
$logs.Cells.Item($newRow, 1).Value = "Opvolging datasheet"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 4).Value = "Documentatie / Datasheets"
$logs.Cells.Item($newRow, 6).Value = "2025-08-28 21:24:45"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Update Dashboard summary: swap "Overig" and "Documentatie / Datasheets" rows,
#     bump "Documentatie / Datasheets" count from 1 to 2 ---
$dashboard.Cells.Item(7, 1).Value = "Documentatie / Datasheets"
$dashboard.Cells.Item(7, 2).Value = 2
$dashboard.Cells.Item(8, 1).Value = "Overig"
$dashboard.Cells.Item(8, 2).Value = 1

# --- Extend the conditional-formatting ranges on Logs so the new row is covered ---
$logs.Range("D2:D29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D30"))
$logs.Range("G2:G29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G30"))
$logs.Range("H2:H29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H30"))
$logs.Range("I2:I29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I30"))
$logs.Range("J2:J29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J30"))
